$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CO")

$ws.Range("D8").Value = 139000
$ws.Range("E8").Value = 112800
$ws.Range("F8").Value = 98400
$ws.Range("G8").Value = 94300
$ws.Range("H8").Value = 85000
$ws.Range("I8").Value = 78100
$ws.Range("J8").Value = 56500
$ws.Range("D9").Value = 26900
$ws.Range("E9").Value = 21200
$ws.Range("F9").Value = 21500
$ws.Range("G9").Value = 19400
$ws.Range("H9").Value = 15800
$ws.Range("I9").Value = 15800
$ws.Range("J9").Value = 12900
$ws.Range("D10").Value = 112100
$ws.Range("E10").Value = 91600
$ws.Range("F10").Value = 76900
$ws.Range("G10").Value = 74900
$ws.Range("H10").Value = 69300
$ws.Range("I10").Value = 62300
$ws.Range("J10").Value = 43600
$ws.Range("F12").Value = 1200
$ws.Range("D15").Value = 2400
$ws.Range("E15").Value = 2400
$ws.Range("I15").Value = 2100
$ws.Range("D17").Value = 97500
$ws.Range("E17").Value = 73900
$ws.Range("F17").Value = 71200
$ws.Range("G17").Value = 59400
$ws.Range("H17").Value = 50600
$ws.Range("I17").Value = 47000
$ws.Range("J17").Value = 36500
$ws.Range("D18").Value = 41500
$ws.Range("E18").Value = 38900
$ws.Range("F18").Value = 27200
$ws.Range("G18").Value = 34900
$ws.Range("H18").Value = 34400
$ws.Range("I18").Value = 31100
$ws.Range("J18").Value = 20000
$ws.Range("D20").Value = 4000
$ws.Range("E20").Value = 3500
$ws.Range("F20").Value = 9800
$ws.Range("G20").Value = 3200
$ws.Range("H20").Value = 4300
$ws.Range("J20").Value = 3000
$ws.Range("D21").Value = 53100
$ws.Range("E21").Value = 49900
$ws.Range("F21").Value = 44500
$ws.Range("G21").Value = 45500
$ws.Range("H21").Value = 44200
$ws.Range("I21").Value = 39200
$ws.Range("J21").Value = 27400
$ws.Range("E22").Value = 17700
$ws.Range("F22").Value = 16000
$ws.Range("G22").Value = 15000
$ws.Range("H22").Value = 10400
$ws.Range("I22").Value = 10400
$ws.Range("D23").Value = 45000
$ws.Range("E23").Value = 24700
$ws.Range("F23").Value = 21000
$ws.Range("G23").Value = 23000
$ws.Range("H23").Value = 28300
$ws.Range("I23").Value = 23500
$ws.Range("J23").Value = 22500
$ws.Range("D24").Value = 9300
$ws.Range("E24").Value = 5600
$ws.Range("F24").Value = 7400
$ws.Range("G24").Value = 7000
$ws.Range("H24").Value = 8700
$ws.Range("I24").Value = 5700
$ws.Range("D26").Value = 35700
$ws.Range("E26").Value = 19100
$ws.Range("F26").Value = 13600
$ws.Range("G26").Value = 16000
$ws.Range("H26").Value = 19700
$ws.Range("I26").Value = 17800
$ws.Range("J26").Value = 21100
$ws.Range("D27").Value = 35200
$ws.Range("E27").Value = 17200
$ws.Range("F27").Value = 13500
$ws.Range("G27").Value = 14700
$ws.Range("H27").Value = 17300
$ws.Range("I27").Value = 16700
$ws.Range("J27").Value = 19600
$ws.Range("D32").Value = -4000
$ws.Range("E32").Value = -3500
$ws.Range("F32").Value = -9800
$ws.Range("G32").Value = -3200
$ws.Range("H32").Value = -4300
$ws.Range("J32").Value = -3000
$ws.Range("D33").Value = 35200
$ws.Range("E33").Value = 17200
$ws.Range("F33").Value = 13500
$ws.Range("G33").Value = 14700
$ws.Range("H33").Value = 17300
$ws.Range("I33").Value = 16700
$ws.Range("J33").Value = 19600
$ws.Range("D35").Value = 35200
$ws.Range("E35").Value = 17200
$ws.Range("F35").Value = 13500
$ws.Range("G35").Value = 14700
$ws.Range("H35").Value = 17300
$ws.Range("I35").Value = 16700
$ws.Range("J35").Value = 19600
$ws.Range("D41").Value = 630800
$ws.Range("E41").Value = 521000
$ws.Range("F41").Value = 446500
$ws.Range("G41").Value = 361600
$ws.Range("H41").Value = 279400
$ws.Range("I41").Value = 221700
$ws.Range("J41").Value = 117900
$ws.Range("D43").Value = 17400
$ws.Range("E43").Value = 17500
$ws.Range("F43").Value = 19600
$ws.Range("G43").Value = 18700
$ws.Range("H43").Value = 16600
$ws.Range("I43").Value = 11200
$ws.Range("J43").Value = 11900
$ws.Range("D44").Value = 4100
$ws.Range("E44").Value = 4600
$ws.Range("F44").Value = 4200
$ws.Range("H44").Value = 4700
$ws.Range("F45").Value = 4700
$ws.Range("G45").Value = 4200
$ws.Range("H45").Value = 4700
$ws.Range("J45").Value = 2300
$ws.Range("D46").Value = 654300
$ws.Range("E46").Value = 544900
$ws.Range("F46").Value = 472800
$ws.Range("G46").Value = 388600
$ws.Range("H46").Value = 305400
$ws.Range("I46").Value = 237200
$ws.Range("J46").Value = 133100
$ws.Range("D47").Value = 66000
$ws.Range("E47").Value = 77900
$ws.Range("F47").Value = 76700
$ws.Range("G47").Value = 75100
$ws.Range("H47").Value = 82900
$ws.Range("I47").Value = 78200
$ws.Range("J47").Value = 72200
$ws.Range("D48").Value = 82100
$ws.Range("E48").Value = 81800
$ws.Range("F48").Value = 85300
$ws.Range("G48").Value = 89500
$ws.Range("H48").Value = 93000
$ws.Range("I48").Value = 69500
$ws.Range("J48").Value = 39800
$ws.Range("D49").Value = 15100
$ws.Range("E49").Value = 15800
$ws.Range("F49").Value = 16500
$ws.Range("G49").Value = 17200
$ws.Range("H49").Value = 17900
$ws.Range("I49").Value = 18600
$ws.Range("J49").Value = 19300
$ws.Range("D52").Value = 49900
$ws.Range("E52").Value = 48700
$ws.Range("F52").Value = 44400
$ws.Range("G52").Value = 40400
$ws.Range("H52").Value = 39600
$ws.Range("I52").Value = 39700
$ws.Range("J52").Value = 6300
$ws.Range("D54").Value = 867400
$ws.Range("E54").Value = 769200
$ws.Range("F54").Value = 695700
$ws.Range("G54").Value = 610200
$ws.Range("H54").Value = 538900
$ws.Range("I54").Value = 443200
$ws.Range("J54").Value = 270700
$ws.Range("F57").Value = 2000
$ws.Range("G57").Value = 1900
$ws.Range("I57").Value = 1500
$ws.Range("J57").Value = 3800
$ws.Range("E58").Value = 153000
$ws.Range("F58").Value = 8900
$ws.Range("G58").Value = 17800
$ws.Range("H58").Value = 8900
$ws.Range("I58").Value = 7400
$ws.Range("J58").Value = 6700
$ws.Range("D59").Value = 67800
$ws.Range("E59").Value = 60100
$ws.Range("F59").Value = 58600
$ws.Range("G59").Value = 51600
$ws.Range("H59").Value = 48500
$ws.Range("I59").Value = 40500
$ws.Range("J59").Value = 18800
$ws.Range("D60").Value = 69500
$ws.Range("E60").Value = 214800
$ws.Range("F60").Value = 67400
$ws.Range("G60").Value = 62400
$ws.Range("H60").Value = 59000
$ws.Range("I60").Value = 49300
$ws.Range("J60").Value = 29300
$ws.Range("F61").Value = 134500
$ws.Range("G61").Value = 119900
$ws.Range("H61").Value = 115400
$ws.Range("I61").Value = 111600
$ws.Range("D62").Value = 335000
$ws.Range("E62").Value = 281000
$ws.Range("F62").Value = 239600
$ws.Range("G62").Value = 198900
$ws.Range("H62").Value = 150800
$ws.Range("I62").Value = 98000
$ws.Range("J62").Value = 58100
$ws.Range("D66").Value = 405300
$ws.Range("E66").Value = 496400
$ws.Range("F66").Value = 442100
$ws.Range("G66").Value = 382000
$ws.Range("H66").Value = 325900
$ws.Range("I66").Value = 259600
$ws.Range("J66").Value = 92200
$ws.Range("D72").Value = 165800
$ws.Range("E72").Value = 130600
$ws.Range("F72").Value = 111800
$ws.Range("G72").Value = 98300
$ws.Range("H72").Value = 82400
$ws.Range("I72").Value = 62800
$ws.Range("J72").Value = 46200
$ws.Range("D76").Value = 462100
$ws.Range("E76").Value = 272800
$ws.Range("F76").Value = 253700
$ws.Range("G76").Value = 228200
$ws.Range("H76").Value = 213000
$ws.Range("I76").Value = 183600
$ws.Range("J76").Value = 178500
$ws.Range("D81").Value = 35200
$ws.Range("E81").Value = 17200
$ws.Range("F81").Value = 13500
$ws.Range("G81").Value = 14700
$ws.Range("H81").Value = 17300
$ws.Range("I81").Value = 16700
$ws.Range("J81").Value = 19600
$ws.Range("D83").Value = 7500
$ws.Range("E83").Value = 7500
$ws.Range("F83").Value = 7400
$ws.Range("G83").Value = 7500
$ws.Range("H83").Value = 5500
$ws.Range("I83").Value = 5300
$ws.Range("J83").Value = 4400
$ws.Range("D89").Value = 121500
$ws.Range("E89").Value = 94600
$ws.Range("F89").Value = 86200
$ws.Range("G89").Value = 88300
$ws.Range("H89").Value = 79500
$ws.Range("I89").Value = 85900
$ws.Range("J89").Value = 52500
$ws.Range("D91").Value = -10000
$ws.Range("E91").Value = -3700
$ws.Range("F91").Value = -2500
$ws.Range("G91").Value = -5600
$ws.Range("H91").Value = -22300
$ws.Range("I91").Value = -33600
$ws.Range("J91").Value = -6400
$ws.Range("D94").Value = -9900
$ws.Range("E94").Value = -13400
$ws.Range("G94").Value = -6300
$ws.Range("H94").Value = -22200
$ws.Range("I94").Value = -73300
$ws.Range("J94").Value = -18200
$ws.Range("E100").Value = -8900
$ws.Range("I100").Value = 91800
$ws.Range("J100").Value = -6600
$ws.Range("D101").Value = -1500
$ws.Range("D102").Value = 109900
$ws.Range("E102").Value = 74500
$ws.Range("F102").Value = 84900
$ws.Range("G102").Value = 82200
$ws.Range("H102").Value = 57700
$ws.Range("I102").Value = 103900
$ws.Range("J102").Value = 27100
